$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Apply the country-reordering + data-refresh edits.
# Each block corresponds to one data row (identified by its row number in
# the "Pais" sheet). Column A gets the (possibly new) country name, columns
# B-H get the refreshed case statistics.

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 646300
$ws.Cells.Item(4, 3).Value = 2211
$ws.Cells.Item(4, 4).Value = 49091
$ws.Cells.Item(4, 5).Value = 568569
$ws.Cells.Item(4, 6).Value = 13487
$ws.Cells.Item(4, 7).Value = 111
$ws.Cells.Item(4, 8).Value = 28640

# Row 15: Brasil -> Brasil
$ws.Cells.Item(15, 1).Value = "Brasil"
$ws.Cells.Item(15, 2).Value = 29165
$ws.Cells.Item(15, 3).Value = 555
$ws.Cells.Item(15, 4).Value = 14026
$ws.Cells.Item(15, 5).Value = 13375
$ws.Cells.Item(15, 6).Value = 296
$ws.Cells.Item(15, 7).Value = 7
$ws.Cells.Item(15, 8).Value = 1764

# Row 47: Panama -> Republica Dominicana
$ws.Cells.Item(47, 1).Value = "Republica Dominicana"
$ws.Cells.Item(47, 2).Value = 3755
$ws.Cells.Item(47, 3).Value = 141
$ws.Cells.Item(47, 4).Value = 215
$ws.Cells.Item(47, 5).Value = 3344
$ws.Cells.Item(47, 6).Value = 143
$ws.Cells.Item(47, 7).Value = 7
$ws.Cells.Item(47, 8).Value = 196

# Row 48: Singapur -> Panama
$ws.Cells.Item(48, 1).Value = "Panama"
$ws.Cells.Item(48, 2).Value = 3751
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 75
$ws.Cells.Item(48, 5).Value = 3573
$ws.Cells.Item(48, 6).Value = 106
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 103

# Row 49: Republica Dominicana -> Singapur
$ws.Cells.Item(49, 1).Value = "Singapur"
$ws.Cells.Item(49, 2).Value = 3699
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 652
$ws.Cells.Item(49, 5).Value = 3037
$ws.Cells.Item(49, 6).Value = 29
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 10

# Row 60: Moldavia -> Moldavia
$ws.Cells.Item(60, 1).Value = "Moldavia"
$ws.Cells.Item(60, 2).Value = 2154
$ws.Cells.Item(60, 3).Value = 105
$ws.Cells.Item(60, 4).Value = 235
$ws.Cells.Item(60, 5).Value = 1866
$ws.Cells.Item(60, 6).Value = 80
$ws.Cells.Item(60, 7).Value = 7
$ws.Cells.Item(60, 8).Value = 53

# Row 77: Oman -> Republica de Macedonia
$ws.Cells.Item(77, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(77, 2).Value = 1081
$ws.Cells.Item(77, 3).Value = 107
$ws.Cells.Item(77, 4).Value = 121
$ws.Cells.Item(77, 5).Value = 914
$ws.Cells.Item(77, 6).Value = 15
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 46

# Row 78: Hong Kong -> Oman
$ws.Cells.Item(78, 1).Value = "Oman"
$ws.Cells.Item(78, 2).Value = 1019
$ws.Cells.Item(78, 3).Value = 109
$ws.Cells.Item(78, 4).Value = 176
$ws.Cells.Item(78, 5).Value = 839
$ws.Cells.Item(78, 6).Value = 3
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 4

# Row 79: Eslovaquia -> Hong Kong
$ws.Cells.Item(79, 1).Value = "Hong Kong"
$ws.Cells.Item(79, 2).Value = 1018
$ws.Cells.Item(79, 3).Value = 1
$ws.Cells.Item(79, 4).Value = 485
$ws.Cells.Item(79, 5).Value = 529
$ws.Cells.Item(79, 6).Value = 9
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 4

# Row 80: Republica de Macedonia -> Eslovaquia
$ws.Cells.Item(80, 1).Value = "Eslovaquia"
$ws.Cells.Item(80, 2).Value = 977
$ws.Cells.Item(80, 3).Value = 114
$ws.Cells.Item(80, 4).Value = 167
$ws.Cells.Item(80, 5).Value = 802
$ws.Cells.Item(80, 6).Value = 5
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 8

# Row 84: Bulgaria -> Bulgaria
$ws.Cells.Item(84, 1).Value = "Bulgaria"
$ws.Cells.Item(84, 2).Value = 800
$ws.Cells.Item(84, 3).Value = 53
$ws.Cells.Item(84, 4).Value = 122
$ws.Cells.Item(84, 5).Value = 640
$ws.Cells.Item(84, 6).Value = 37
$ws.Cells.Item(84, 7).Value = 2
$ws.Cells.Item(84, 8).Value = 38

# Row 100: Republica de Yibuti -> Guinea
$ws.Cells.Item(100, 1).Value = "Guinea"
$ws.Cells.Item(100, 2).Value = 438
$ws.Cells.Item(100, 3).Value = 34
$ws.Cells.Item(100, 4).Value = 49
$ws.Cells.Item(100, 5).Value = 388
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 1

# Row 101: Honduras -> Republica de Yibuti
$ws.Cells.Item(101, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(101, 2).Value = 435
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 71
$ws.Cells.Item(101, 5).Value = 362
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 2

# Row 102: San Marino -> Honduras
$ws.Cells.Item(102, 1).Value = "Honduras"
$ws.Cells.Item(102, 2).Value = 426
$ws.Cells.Item(102, 3).Value = 7
$ws.Cells.Item(102, 4).Value = 9
$ws.Cells.Item(102, 5).Value = 382
$ws.Cells.Item(102, 6).Value = 10
$ws.Cells.Item(102, 7).Value = 4
$ws.Cells.Item(102, 8).Value = 35

# Row 103: Malta -> San Marino
$ws.Cells.Item(103, 1).Value = "San Marino"
$ws.Cells.Item(103, 2).Value = 426
$ws.Cells.Item(103, 3).Value = 33
$ws.Cells.Item(103, 4).Value = 55
$ws.Cells.Item(103, 5).Value = 333
$ws.Cells.Item(103, 6).Value = 15
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 38

# Row 104: Nigeria -> Malta
$ws.Cells.Item(104, 1).Value = "Malta"
$ws.Cells.Item(104, 2).Value = 412
$ws.Cells.Item(104, 3).Value = 13
$ws.Cells.Item(104, 4).Value = 82
$ws.Cells.Item(104, 5).Value = 327
$ws.Cells.Item(104, 6).Value = 4
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 3

# Row 105: Guinea -> Nigeria
$ws.Cells.Item(105, 1).Value = "Nigeria"
$ws.Cells.Item(105, 2).Value = 407
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 128
$ws.Cells.Item(105, 5).Value = 267
$ws.Cells.Item(105, 6).Value = 2
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 12

# Row 169: Mongolia -> Mozambique
$ws.Cells.Item(169, 1).Value = "Mozambique"
$ws.Cells.Item(169, 2).Value = 31
$ws.Cells.Item(169, 3).Value = 2
$ws.Cells.Item(169, 4).Value = 2
$ws.Cells.Item(169, 5).Value = 29
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

# Row 170: Mozambique -> Mongolia
$ws.Cells.Item(170, 1).Value = "Mongolia"
$ws.Cells.Item(170, 2).Value = 31
$ws.Cells.Item(170, 3).Value = 1
$ws.Cells.Item(170, 4).Value = 5
$ws.Cells.Item(170, 5).Value = 26
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 0

# Row 180: Fiyi -> Islas Virgenes de los Estados Unidos
$ws.Cells.Item(180, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(180, 2).Value = 17
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 17
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

# Row 181: Islas Virgenes de los Estados Unidos -> Fiyi
$ws.Cells.Item(181, 1).Value = "Fiyi"
$ws.Cells.Item(181, 2).Value = 17
$ws.Cells.Item(181, 3).Value = 1
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 17
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

# Row 185: Suazilandia -> Dominica
$ws.Cells.Item(185, 1).Value = "Dominica"
$ws.Cells.Item(185, 2).Value = 16
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 8
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 186: Dominica -> Suazilandia
$ws.Cells.Item(186, 1).Value = "Suazilandia"
$ws.Cells.Item(186, 2).Value = 16
$ws.Cells.Item(186, 3).Value = 1
$ws.Cells.Item(186, 4).Value = 8
$ws.Cells.Item(186, 5).Value = 8
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

# Row 190: San Cristobal y Nieves -> Granada
$ws.Cells.Item(190, 1).Value = "Granada"
$ws.Cells.Item(190, 2).Value = 14
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 14
$ws.Cells.Item(190, 6).Value = 2
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 191: Granada -> San Cristobal y Nieves
$ws.Cells.Item(191, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(191, 2).Value = 14
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 14
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0

# Row 209: Sudan del Sur -> Santo Tome y Principe
$ws.Cells.Item(209, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(209, 2).Value = 4
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 0
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 210: Santo Tome y Principe -> Sudan del Sur
$ws.Cells.Item(210, 1).Value = "Sudan del Sur"
$ws.Cells.Item(210, 2).Value = 4
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 0
$ws.Cells.Item(210, 5).Value = 4
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0
